$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 13153.818
$ws.Range("J28").Value = 20141.143
$ws.Range("L28").Value = 20141.143
$ws.Range("N28").Value = -21111.143
$ws.Range("H40").Value = 5228.8887
$ws.Range("I40").Value = 4119.909
$ws.Range("J40").Value = 6971.5713
$ws.Range("K40").Value = 4119.909
$ws.Range("L40").Value = 6971.5713
$ws.Range("M40").Value = -3944.909
$ws.Range("N40").Value = -7321.5713
$ws.Range("H41").Value = 2448.7144
$ws.Range("I41").Value = 2298
$ws.Range("J41").Value = 2649.6667
$ws.Range("K41").Value = 2298
$ws.Range("L41").Value = 2649.6667
$ws.Range("M41").Value = -1858
$ws.Range("N41").Value = -3529.6667
$ws.Range("H53").Value = 964.625
$ws.Range("I53").Value = 548.5454999999999
$ws.Range("J53").Value = 1880
$ws.Range("K53").Value = 548.5454999999999
$ws.Range("L53").Value = 1880
$ws.Range("M53").Value = 88.45450000000005
$ws.Range("N53").Value = -3154
$ws.Range("H62").Value = 8267.214
$ws.Range("I62").Value = 4960.5
$ws.Range("J62").Value = 10747.25
$ws.Range("K62").Value = 4960.5
$ws.Range("L62").Value = 10747.25
$ws.Range("M62").Value = -4336.5
$ws.Range("N62").Value = -11995.25
$ws.Range("H65").Value = 8267.214
$ws.Range("I65").Value = 4960.5
$ws.Range("J65").Value = 10747.25
$ws.Range("K65").Value = 24802.5
$ws.Range("L65").Value = 53736.25
$ws.Range("M65").Value = -21682.5
$ws.Range("N65").Value = -59976.25
$ws.Range("H76").Value = 3650
$ws.Range("J76").Value = 3999.6667
$ws.Range("L76").Value = 3999.6667
$ws.Range("N76").Value = -4629.6667
$ws.Range("H79").Value = 3650
$ws.Range("J79").Value = 3999.6667
$ws.Range("L79").Value = 3999.6667
$ws.Range("N79").Value = -6183.6667
$ws.Range("H86").Value = 1499.3334
$ws.Range("I86").Value = 1249.5
$ws.Range("K86").Value = 1249.5
$ws.Range("M86").Value = -126.5
$ws.Range("H89").Value = 1499.3334
$ws.Range("I89").Value = 1249.5
$ws.Range("K89").Value = 6247.5
$ws.Range("M89").Value = -631.5
$ws.Range("H98").Value = 1701.625
$ws.Range("I98").Value = 1587.5714
$ws.Range("K98").Value = 1587.5714
$ws.Range("M98").Value = -89.57140000000004
$ws.Range("H106").Value = 6200.375
$ws.Range("I106").Value = 6200.375
$ws.Range("K106").Value = 6200.375
$ws.Range("M106").Value = -5569.375
$ws.Range("H112").Value = 2203.25
$ws.Range("I112").Value = 1997.5
$ws.Range("J112").Value = 2409
$ws.Range("K112").Value = 5992.5
$ws.Range("L112").Value = 7227
$ws.Range("M112").Value = -4884.5
$ws.Range("N112").Value = -9443
$ws.Range("H122").Value = 1701.625
$ws.Range("I122").Value = 1587.5714
$ws.Range("K122").Value = 4762.7142
$ws.Range("M122").Value = -2312.7142
$ws.Range("H135").Value = 1915.6666
$ws.Range("I135").Value = 1558.2858
$ws.Range("K135").Value = 14024.5722
$ws.Range("M135").Value = -11489.5722
$ws.Range("H137").Value = 2177.92
$ws.Range("I137").Value = 1378.4375
$ws.Range("J137").Value = 3599.2222
$ws.Range("K137").Value = 4135.3125
$ws.Range("L137").Value = 10797.6666
$ws.Range("M137").Value = -1585.3125
$ws.Range("N137").Value = -15897.6666
$ws.Range("H140").Value = 72500
$ws.Range("J140").Value = 125000
$ws.Range("L140").Value = 125000
$ws.Range("N140").Value = -135360

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 3871.875
$ws.Range("I45").Value = 2795
$ws.Range("K45").Value = 2795
$ws.Range("M45").Value = -2418
$ws.Range("H74").Value = 7523.643
$ws.Range("I74").Value = 7170.636
$ws.Range("J74").Value = 8818
$ws.Range("K74").Value = 7170.636
$ws.Range("L74").Value = 8818
$ws.Range("M74").Value = -6296.636
$ws.Range("N74").Value = -10566
$ws.Range("H77").Value = 7523.643
$ws.Range("I77").Value = 7170.636
$ws.Range("J77").Value = 8818
$ws.Range("K77").Value = 35853.18
$ws.Range("L77").Value = 44090
$ws.Range("M77").Value = -31485.18
$ws.Range("N77").Value = -52826
$ws.Range("H97").Value = 2311.1667
$ws.Range("I97").Value = 1765
$ws.Range("J97").Value = 2857.3333
$ws.Range("K97").Value = 1765
$ws.Range("L97").Value = 2857.3333
$ws.Range("M97").Value = -1269
$ws.Range("N97").Value = -3849.3333
$ws.Range("H132").Value = 3085.6155
$ws.Range("I132").Value = 1790.4445
$ws.Range("J132").Value = 5999.75
$ws.Range("K132").Value = 5371.333500000001
$ws.Range("L132").Value = 17999.25
$ws.Range("M132").Value = -2841.333500000001
$ws.Range("N132").Value = -23059.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3524.8667
$ws.Range("I20").Value = 3387.7273
$ws.Range("J20").Value = 3902
$ws.Range("K20").Value = 3387.7273
$ws.Range("L20").Value = 3902
$ws.Range("M20").Value = -3140.7273
$ws.Range("N20").Value = -4396
$ws.Range("H105").Value = 2180.9
$ws.Range("I105").Value = 2102.375
$ws.Range("J105").Value = 2495
$ws.Range("K105").Value = 2102.375
$ws.Range("L105").Value = 2495
$ws.Range("M105").Value = -355.375
$ws.Range("N105").Value = -5989

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1892
$ws.Range("I22").Value = 1923.6
$ws.Range("K22").Value = 1923.6
$ws.Range("M22").Value = -1573.6
$ws.Range("H99").Value = 3632.85
$ws.Range("I99").Value = 3172.5625
$ws.Range("K99").Value = 3172.5625
$ws.Range("M99").Value = -1674.5625
$ws.Range("H126").Value = 3632.85
$ws.Range("I126").Value = 3172.5625
$ws.Range("K126").Value = 9517.6875
$ws.Range("M126").Value = -7047.6875

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H139").Value = 6272
$ws.Range("I139").Value = 5974.2
$ws.Range("J139").Value = 7016.5
$ws.Range("K139").Value = 17922.6
$ws.Range("L139").Value = 21049.5
$ws.Range("M139").Value = -12782.6
$ws.Range("N139").Value = -31329.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H92").Value = 6083.1665
$ws.Range("I92").Value = 0
$ws.Range("J92").Value = 6083.1665
$ws.Range("K92").Value = 0
$ws.Range("M92").ClearContents()
$ws.Range("N92").Value = -9827.166499999999
$ws.Range("H102").Value = 2994.25
$ws.Range("I102").Value = 2994.25
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 2994.25
$ws.Range("L102").Value = 0
$ws.Range("M102").ClearContents()
$ws.Range("N102").Value = -1372.25
$ws.Range("H126").Value = 5666.6665
$ws.Range("I126").Value = 3500
$ws.Range("J126").Value = 6750
$ws.Range("K126").Value = 10500
$ws.Range("L126").Value = 20250
$ws.Range("M126").Value = -8030
$ws.Range("N126").Value = -25190
$ws.Range("H132").Value = 43051.777
$ws.Range("I132").Value = 60914.223
$ws.Range("J132").Value = 7326.8887
$ws.Range("K132").Value = 182742.669
$ws.Range("L132").Value = 21980.6661
$ws.Range("M132").Value = -180212.669
$ws.Range("N132").Value = -27040.6661

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H5").Value = 63010
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = 63010
$ws.Range("K5").Value = 0
$ws.Range("M5").ClearContents()
$ws.Range("N5").Value = -63236
$ws.Range("H7").Value = 7677.5625
$ws.Range("I7").Value = 7211.8335
$ws.Range("J7").Value = 7957
$ws.Range("K7").Value = 7211.8335
$ws.Range("L7").Value = 7957
$ws.Range("M7").Value = -7099.8335
$ws.Range("N7").Value = -8181
$ws.Range("H22").Value = 1200
$ws.Range("I22").Value = 400
$ws.Range("J22").Value = 2000
$ws.Range("K22").Value = 400
$ws.Range("L22").Value = 2000
$ws.Range("M22").Value = -105
$ws.Range("N22").Value = -2590
$ws.Range("H27").Value = 1200
$ws.Range("I27").Value = 400
$ws.Range("J27").Value = 2000
$ws.Range("K27").Value = 400
$ws.Range("L27").Value = 2000
$ws.Range("M27").Value = -293
$ws.Range("N27").Value = -2214
$ws.Range("H63").Value = 42959.332
$ws.Range("J63").Value = 40000
$ws.Range("L63").Value = 40000
$ws.Range("N63").Value = -41498
$ws.Range("H66").Value = 42959.332
$ws.Range("J66").Value = 40000
$ws.Range("L66").Value = 120000
$ws.Range("N66").Value = -127488
$ws.Range("H126").Value = 7677.5625
$ws.Range("I126").Value = 7211.8335
$ws.Range("J126").Value = 7957
$ws.Range("K126").Value = 21635.5005
$ws.Range("L126").Value = 23871
$ws.Range("M126").Value = -19165.5005
$ws.Range("N126").Value = -28811
$ws.Range("H130").Value = 25330
$ws.Range("J130").Value = 25330
$ws.Range("L130").Value = 25330
$ws.Range("N130").Value = -35370
$ws.Range("H132").Value = 3602.2307
$ws.Range("I132").Value = 3610.8333
$ws.Range("J132").Value = 3499
$ws.Range("K132").Value = 10832.4999
$ws.Range("L132").Value = 10497
$ws.Range("M132").Value = -8302.499899999999
$ws.Range("N132").Value = -15557
$ws.Range("H136").Value = 3938.3333
$ws.Range("I136").Value = 2542.75
$ws.Range("J136").Value = 6729.5
$ws.Range("K136").Value = 7628.25
$ws.Range("L136").Value = 20188.5
$ws.Range("M136").Value = -5078.25
$ws.Range("N136").Value = -25288.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H124").Value = 0
$ws.Range("J124").Value = 0
$ws.Range("N124").ClearContents()
$ws.Range("H136").Value = 2385.2354
$ws.Range("I136").Value = 1287
$ws.Range("J136").Value = 4398.6665
$ws.Range("K136").Value = 3861
$ws.Range("L136").Value = 13195.9995
$ws.Range("M136").Value = -1311
$ws.Range("N136").Value = -18295.9995
